# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit data refresh to the Sheets workbook
# (columns H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ)

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (ALC) ---
$ws = $wb.Worksheets.Item(1)
# Row 74
$ws.Cells.Item(74, 8).Value = 0  # H74
$ws.Cells.Item(74, 9).Value = 0  # I74
$ws.Cells.Item(74, 11).Value = 0  # K74
$ws.Cells.Item(74, 13).ClearContents()  # M74
# Row 77
$ws.Cells.Item(77, 8).Value = 0  # H77
$ws.Cells.Item(77, 9).Value = 0  # I77
$ws.Cells.Item(77, 11).Value = 0  # K77
$ws.Cells.Item(77, 13).ClearContents()  # M77
# Row 113
$ws.Cells.Item(113, 8).Value = 0  # H113
$ws.Cells.Item(113, 9).Value = 0  # I113
$ws.Cells.Item(113, 10).Value = 0  # J113
$ws.Cells.Item(113, 11).Value = 0  # K113
$ws.Cells.Item(113, 12).Value = 0  # L113
$ws.Cells.Item(113, 13).ClearContents()  # M113
$ws.Cells.Item(113, 14).ClearContents()  # N113
# Row 127
$ws.Cells.Item(127, 8).Value = 4497.6  # H127
$ws.Cells.Item(127, 10).Value = 4497.6  # J127
$ws.Cells.Item(127, 12).Value = 13492.8  # L127
$ws.Cells.Item(127, 14).Value = -23412.8  # N127

# --- Sheet 2 (ARM) ---
$ws = $wb.Worksheets.Item(2)
# Row 2
$ws.Cells.Item(2, 8).Value = 3058.25  # H2
$ws.Cells.Item(2, 9).Value = 3058.25  # I2
$ws.Cells.Item(2, 10).Value = 0  # J2
$ws.Cells.Item(2, 11).Value = 3058.25  # K2
$ws.Cells.Item(2, 12).Value = 0  # L2
$ws.Cells.Item(2, 13).Value = -2945.25  # M2
$ws.Cells.Item(2, 14).ClearContents()  # N2
# Row 63
$ws.Cells.Item(63, 8).Value = 7908.5  # H63
$ws.Cells.Item(63, 9).Value = 8113.75  # I63
$ws.Cells.Item(63, 10).Value = 7498  # J63
$ws.Cells.Item(63, 11).Value = 8113.75  # K63
$ws.Cells.Item(63, 12).Value = 7498  # L63
$ws.Cells.Item(63, 13).Value = -7427.75  # M63
$ws.Cells.Item(63, 14).Value = -8870  # N63
# Row 66
$ws.Cells.Item(66, 8).Value = 7908.5  # H66
$ws.Cells.Item(66, 9).Value = 8113.75  # I66
$ws.Cells.Item(66, 10).Value = 7498  # J66
$ws.Cells.Item(66, 11).Value = 40568.75  # K66
$ws.Cells.Item(66, 12).Value = 37490  # L66
$ws.Cells.Item(66, 13).Value = -37136.75  # M66
$ws.Cells.Item(66, 14).Value = -44354  # N66
# Row 116
$ws.Cells.Item(116, 8).Value = 3058.25  # H116
$ws.Cells.Item(116, 9).Value = 3058.25  # I116
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 11).Value = 3058.25  # K116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 13).Value = -764.25  # M116
$ws.Cells.Item(116, 14).ClearContents()  # N116

# --- Sheet 3 (BSM) ---
$ws = $wb.Worksheets.Item(3)
# Row 3
$ws.Cells.Item(3, 8).Value = 3058.25  # H3
$ws.Cells.Item(3, 9).Value = 3058.25  # I3
$ws.Cells.Item(3, 10).Value = 0  # J3
$ws.Cells.Item(3, 11).Value = 3058.25  # K3
$ws.Cells.Item(3, 12).Value = 0  # L3
$ws.Cells.Item(3, 13).Value = -2944.25  # M3
$ws.Cells.Item(3, 14).ClearContents()  # N3
# Row 80
$ws.Cells.Item(80, 8).Value = 2147.6667  # H80
$ws.Cells.Item(80, 9).Value = 2471.5  # I80
$ws.Cells.Item(80, 10).Value = 1500  # J80
$ws.Cells.Item(80, 11).Value = 2471.5  # K80
$ws.Cells.Item(80, 12).Value = 1500  # L80
$ws.Cells.Item(80, 13).Value = -1473.5  # M80
$ws.Cells.Item(80, 14).Value = -3496  # N80
# Row 83
$ws.Cells.Item(83, 8).Value = 2147.6667  # H83
$ws.Cells.Item(83, 9).Value = 2471.5  # I83
$ws.Cells.Item(83, 10).Value = 1500  # J83
$ws.Cells.Item(83, 11).Value = 12357.5  # K83
$ws.Cells.Item(83, 12).Value = 7500  # L83
$ws.Cells.Item(83, 13).Value = -7365.5  # M83
$ws.Cells.Item(83, 14).Value = -17484  # N83
# Row 86
$ws.Cells.Item(86, 8).Value = 4000  # H86
$ws.Cells.Item(86, 9).Value = 3000  # I86
$ws.Cells.Item(86, 11).Value = 3000  # K86
$ws.Cells.Item(86, 13).Value = -1877  # M86
# Row 89
$ws.Cells.Item(89, 8).Value = 4000  # H89
$ws.Cells.Item(89, 9).Value = 3000  # I89
$ws.Cells.Item(89, 11).Value = 15000  # K89
$ws.Cells.Item(89, 13).Value = -9384  # M89

# --- Sheet 4 (CRP) ---
$ws = $wb.Worksheets.Item(4)
# Row 22
$ws.Cells.Item(22, 8).Value = 298.83334  # H22
$ws.Cells.Item(22, 9).Value = 298.83334  # I22
$ws.Cells.Item(22, 11).Value = 298.83334  # K22
$ws.Cells.Item(22, 13).Value = 51.16665999999998  # M22
# Row 32
$ws.Cells.Item(32, 8).Value = 500  # H32
$ws.Cells.Item(32, 9).Value = 500  # I32
$ws.Cells.Item(32, 11).Value = 500  # K32
$ws.Cells.Item(32, 13).Value = -184  # M32

# --- Sheet 5 (CUL) ---
$ws = $wb.Worksheets.Item(5)
# Row 8
$ws.Cells.Item(8, 8).Value = 916  # H8
$ws.Cells.Item(8, 9).Value = 916  # I8
$ws.Cells.Item(8, 11).Value = 2748  # K8
$ws.Cells.Item(8, 13).Value = -2609  # M8
# Row 23
$ws.Cells.Item(23, 8).Value = 0  # H23
$ws.Cells.Item(23, 9).Value = 0  # I23
$ws.Cells.Item(23, 11).Value = 0  # K23
$ws.Cells.Item(23, 13).ClearContents()  # M23
# Row 33
$ws.Cells.Item(33, 8).Value = 492.25  # H33
$ws.Cells.Item(33, 9).Value = 643  # I33
$ws.Cells.Item(33, 10).Value = 40  # J33
$ws.Cells.Item(33, 11).Value = 3858  # K33
$ws.Cells.Item(33, 12).Value = 240  # L33
$ws.Cells.Item(33, 13).Value = -3575  # M33
$ws.Cells.Item(33, 14).Value = -806  # N33
# Row 38
$ws.Cells.Item(38, 8).Value = 87.5  # H38
$ws.Cells.Item(38, 9).Value = 0  # I38
$ws.Cells.Item(38, 10).Value = 87.5  # J38
$ws.Cells.Item(38, 11).Value = 0  # K38
$ws.Cells.Item(38, 12).Value = 262.5  # L38
$ws.Cells.Item(38, 13).ClearContents()  # M38
$ws.Cells.Item(38, 14).Value = -956.5  # N38
# Row 44
$ws.Cells.Item(44, 8).Value = 498  # H44
$ws.Cells.Item(44, 9).Value = 498  # I44
$ws.Cells.Item(44, 10).Value = 0  # J44
$ws.Cells.Item(44, 11).Value = 1494  # K44
$ws.Cells.Item(44, 12).Value = 0  # L44
$ws.Cells.Item(44, 13).Value = -1096  # M44
$ws.Cells.Item(44, 14).ClearContents()  # N44
# Row 47
$ws.Cells.Item(47, 8).Value = 703  # H47
$ws.Cells.Item(47, 9).Value = 703  # I47
$ws.Cells.Item(47, 11).Value = 2109  # K47
$ws.Cells.Item(47, 13).Value = -1678  # M47
# Row 107
$ws.Cells.Item(107, 8).Value = 741.7857  # H107
$ws.Cells.Item(107, 10).Value = 798.3333  # J107
$ws.Cells.Item(107, 12).Value = 2394.9999  # L107
$ws.Cells.Item(107, 14).Value = -6234.9999  # N107

# --- Sheet 6 (GSM) ---
$ws = $wb.Worksheets.Item(6)
# Row 31
$ws.Cells.Item(31, 8).Value = 1149.5  # H31
$ws.Cells.Item(31, 9).Value = 1149.5  # I31
$ws.Cells.Item(31, 11).Value = 1149.5  # K31
$ws.Cells.Item(31, 13).Value = -857.5  # M31
# Row 37
$ws.Cells.Item(37, 8).Value = 1149.5  # H37
$ws.Cells.Item(37, 9).Value = 1149.5  # I37
$ws.Cells.Item(37, 11).Value = 1149.5  # K37
$ws.Cells.Item(37, 13).Value = -872.5  # M37
# Row 55
$ws.Cells.Item(55, 8).Value = 34993  # H55
$ws.Cells.Item(55, 9).Value = 0  # I55
$ws.Cells.Item(55, 10).Value = 34993  # J55
$ws.Cells.Item(55, 11).Value = 0  # K55
$ws.Cells.Item(55, 12).Value = 34993  # L55
$ws.Cells.Item(55, 13).ClearContents()  # M55
$ws.Cells.Item(55, 14).Value = -35647  # N55
# Row 64
$ws.Cells.Item(64, 8).Value = 0  # H64
$ws.Cells.Item(64, 10).Value = 0  # J64
$ws.Cells.Item(64, 12).Value = 0  # L64
$ws.Cells.Item(64, 14).ClearContents()  # N64
# Row 67
$ws.Cells.Item(67, 8).Value = 0  # H67
$ws.Cells.Item(67, 10).Value = 0  # J67
$ws.Cells.Item(67, 12).Value = 0  # L67
$ws.Cells.Item(67, 14).ClearContents()  # N67

# --- Sheet 7 (LTW) ---
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Cells.Item(7, 8).Value = 0  # H7
$ws.Cells.Item(7, 9).Value = 0  # I7
$ws.Cells.Item(7, 11).Value = 0  # K7
$ws.Cells.Item(7, 13).ClearContents()  # M7
# Row 55
$ws.Cells.Item(55, 8).Value = 1708.3334  # H55
$ws.Cells.Item(55, 9).Value = 849.25  # I55
$ws.Cells.Item(55, 10).Value = 2395.6  # J55
$ws.Cells.Item(55, 11).Value = 849.25  # K55
$ws.Cells.Item(55, 12).Value = 2395.6  # L55
$ws.Cells.Item(55, 13).Value = -676.25  # M55
$ws.Cells.Item(55, 14).Value = -2741.6  # N55
# Row 61
$ws.Cells.Item(61, 8).Value = 2000  # H61
$ws.Cells.Item(61, 9).Value = 0  # I61
$ws.Cells.Item(61, 11).Value = 0  # K61
$ws.Cells.Item(61, 13).ClearContents()  # M61
# Row 62
$ws.Cells.Item(62, 8).Value = 35000  # H62
$ws.Cells.Item(62, 10).Value = 35000  # J62
$ws.Cells.Item(62, 12).Value = 35000  # L62
$ws.Cells.Item(62, 14).Value = -36248  # N62
# Row 65
$ws.Cells.Item(65, 8).Value = 35000  # H65
$ws.Cells.Item(65, 10).Value = 35000  # J65
$ws.Cells.Item(65, 12).Value = 105000  # L65
$ws.Cells.Item(65, 14).Value = -111240  # N65
# Row 68
$ws.Cells.Item(68, 8).Value = 3000  # H68
$ws.Cells.Item(68, 9).Value = 3000  # I68
$ws.Cells.Item(68, 11).Value = 3000  # K68
$ws.Cells.Item(68, 13).Value = -2251  # M68
# Row 71
$ws.Cells.Item(71, 8).Value = 3000  # H71
$ws.Cells.Item(71, 9).Value = 3000  # I71
$ws.Cells.Item(71, 11).Value = 15000  # K71
$ws.Cells.Item(71, 13).Value = -11256  # M71
# Row 113
$ws.Cells.Item(113, 8).Value = 2000  # H113
$ws.Cells.Item(113, 9).Value = 0  # I113
$ws.Cells.Item(113, 11).Value = 0  # K113
$ws.Cells.Item(113, 13).ClearContents()  # M113
# Row 126
$ws.Cells.Item(126, 8).Value = 0  # H126
$ws.Cells.Item(126, 9).Value = 0  # I126
$ws.Cells.Item(126, 11).Value = 0  # K126
$ws.Cells.Item(126, 13).ClearContents()  # M126
# Row 136
$ws.Cells.Item(136, 8).Value = 8850.666999999999  # H136
$ws.Cells.Item(136, 10).Value = 6000  # J136
$ws.Cells.Item(136, 12).Value = 18000  # L136
$ws.Cells.Item(136, 14).Value = -23100  # N136

# --- Sheet 8 (WVR) ---
$ws = $wb.Worksheets.Item(8)
# Row 126
$ws.Cells.Item(126, 8).Value = 1518.375  # H126
$ws.Cells.Item(126, 9).Value = 1163.8572  # I126
$ws.Cells.Item(126, 11).Value = 3491.5716  # K126
$ws.Cells.Item(126, 13).Value = -1021.5716  # M126
